$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts B:F left to A:E,
# removing the stray duplicate data that was in column A.
$ws.Range("A:A").Delete()
